$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 37
$ws.Range("C35").Value = 10
$ws.Range("D35").Value = 20
$ws.Range("E35").Value = 17
$ws.Range("F35").Value = 67
$ws.Range("G35").Value = 84

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 42
$ws.Range("C36").Value = 11
$ws.Range("D36").Value = 14
$ws.Range("E36").Value = 24
$ws.Range("F36").Value = 67
$ws.Range("G36").Value = 91
